$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the header of column B from "commentaire" to "nom"
$ws.Range("B1").Value = "nom"

# Restore the active selection to B2 (matches the committed file's sheet view state)
$ws.Range("B2").Select()
